$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the GADM1 NUTS2 SCI Data description cell (A4) to quote the script path.
$ws.Range("A4").Value = "GADM1 NUTS2 SCI Data (Note: in the script ""_prep_scripts/1_setup_map_and_SCI.R"", we filter to create  ""_raw_data/SCI_Nuts2_Nuts2.csv"", which we will use for the rest of the paper)"

# Update the "NUTS 2016 Shapefiles" description cell (A2) to add a note about
# where the shapefiles are downloaded.
$ws.Range("A2").Value = "NUTS 2016 Shapefiles (Note: these are downloaded in the script ""_prep_scripts/1_setup_map_and_SCI.R"")"

# Move the active selection from A4 to A5.
$ws.Range("A5").Select()
